$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing row 259 values (B and D columns were recomputed upstream)
$ws.Cells.Item(259, 2).Value = 6243090940000
$ws.Cells.Item(259, 4).Value = 286710429898.6448

# Append new rows 260-262 with the same structure/style as the existing data rows.
# Copy the formatting of column A's date cell (A259) down into the new rows so the
# new dates pick up the same style (border/font/number-format/alignment) rather
# than minting brand-new style entries.
$ws.Range("A259").Copy() | Out-Null
$ws.Range("A260:A262").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @{ Row = 260; A = 45108; B = 6355692770000;   C = 0.04603977376055174; D = 292614657222.3744 },
    @{ Row = 261; A = 45139; B = 6337051350000;   C = 0.04504991530615922; D = 285483626608.2819 },
    @{ Row = 262; A = 45170; B = 6359425540000;   C = 0.04329632219391124; D = 275339737148.028 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
}
